$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-06-27 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-28 Friday", 2) | Out-Null
$d.Content.Find.Execute("86×48=4128", $true, $false, $false, $false, $false, $true, 1, $false, "37×90=3330", 2) | Out-Null
$d.Content.Find.Execute("25×83=2075", $true, $false, $false, $false, $false, $true, 1, $false, "47×60=2820", 2) | Out-Null
$d.Content.Find.Execute("86×54=4644", $true, $false, $false, $false, $false, $true, 1, $false, "84×90=7560", 2) | Out-Null
$d.Content.Find.Execute("91×66=6006", $true, $false, $false, $false, $false, $true, 1, $false, "35×72=2520", 2) | Out-Null
$d.Content.Find.Execute("84×24=2016", $true, $false, $false, $false, $false, $true, 1, $false, "22×89=1958", 2) | Out-Null
$d.Content.Find.Execute("18×63=1134", $true, $false, $false, $false, $false, $true, 1, $false, "42×36=1512", 2) | Out-Null
$d.Content.Find.Execute("81×90=7290", $true, $false, $false, $false, $false, $true, 1, $false, "13×23=299", 2) | Out-Null
$d.Content.Find.Execute("23×56=1288", $true, $false, $false, $false, $false, $true, 1, $false, "21×79=1659", 2) | Out-Null
$d.Content.Find.Execute("68×65=4420", $true, $false, $false, $false, $false, $true, 1, $false, "93×91=8463", 2) | Out-Null
$d.Content.Find.Execute("68×41=2788", $true, $false, $false, $false, $false, $true, 1, $false, "74×26=1924", 2) | Out-Null
$d.Content.Find.Execute("88×93=8184", $true, $false, $false, $false, $false, $true, 1, $false, "86×91=7826", 2) | Out-Null
$d.Content.Find.Execute("80×86=6880", $true, $false, $false, $false, $false, $true, 1, $false, "61×98=5978", 2) | Out-Null
$d.Content.Find.Execute("53×41=2173", $true, $false, $false, $false, $false, $true, 1, $false, "36×97=3492", 2) | Out-Null
$d.Content.Find.Execute("51×58=2958", $true, $false, $false, $false, $false, $true, 1, $false, "19×94=1786", 2) | Out-Null
$d.Content.Find.Execute("26×83=2158", $true, $false, $false, $false, $false, $true, 1, $false, "75×75=5625", 2) | Out-Null
$d.Content.Find.Execute("75×80=6000", $true, $false, $false, $false, $false, $true, 1, $false, "89×97=8633", 2) | Out-Null
$d.Content.Find.Execute("47×97=4559", $true, $false, $false, $false, $false, $true, 1, $false, "78×39=3042", 2) | Out-Null
$d.Content.Find.Execute("28×33=924", $true, $false, $false, $false, $false, $true, 1, $false, "74×52=3848", 2) | Out-Null
$d.Content.Find.Execute("29×77=2233", $true, $false, $false, $false, $false, $true, 1, $false, "57×16=912", 2) | Out-Null
$d.Content.Find.Execute("32×28=896", $true, $false, $false, $false, $false, $true, 1, $false, "59×55=3245", 2) | Out-Null
$d.Content.Find.Execute("51×11=561", $true, $false, $false, $false, $false, $true, 1, $false, "48×16=768", 2) | Out-Null
$d.Content.Find.Execute("83×49=4067", $true, $false, $false, $false, $false, $true, 1, $false, "19×74=1406", 2) | Out-Null
$d.Content.Find.Execute("20×74=1480", $true, $false, $false, $false, $false, $true, 1, $false, "19×71=1349", 2) | Out-Null
$d.Content.Find.Execute("93×90=8370", $true, $false, $false, $false, $false, $true, 1, $false, "87×68=5916", 2) | Out-Null
$d.Content.Find.Execute("96×98=9408", $true, $false, $false, $false, $false, $true, 1, $false, "81×42=3402", 2) | Out-Null

Write-Output "Replacements applied: 26"
